# "Removed Test Case Inter-Dependency"
#
# The ProductLoanInput/ProductLoanOutput test data reused the same
# product name/short name as other loan-product test cases, causing
# tests to clash when run together. Give this product its own unique
# name and short name, and make the input sheet (where the values are
# entered) the active tab instead of the output sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# productname (shown on both the input sheet and the output/verify sheet)
$ws1.Range("B1").Value = "4135-RBI-EI-FL-SAR-NOREC-MOREREPAY-1st"
$ws2.Range("B1").Value = "4135-RBI-EI-FL-SAR-NOREC-MOREREPAY-1st"

# shortname - was the bare numeric 4135, now a distinct text code
$ws1.Range("B2").Value = "413u"

# Make the input sheet the active/selected tab
$ws1.Activate()
